# Update terminal storage amounts in "input_param" files
# Insert two new input columns ("terminal compressed hydrogen storage amount
# (days)" and "terminal liquid hydrogen storage amount (days)") right after
# the existing "terminal formic acid storage amount (days)" column (AA),
# pushing the dehydrogenation-related columns (old AB..AL) two slots to the
# right (new AD..AN).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two blank columns before the old "AB" column -----------------
$ws.Range("AB1:AC1").EntireColumn.Insert()

# --- 2. Populate the two new header cells (row 1) ---------------------------
$ws.Range("AB1").Value = "terminal compressed hydrogen storage amount (days)"
$ws.Range("AC1").Value = "terminal liquid hydrogen storage amount (days)"

# --- 3. Populate the new data columns for the three data rows ---------------
$ws.Range("AB2").Value = 0.25
$ws.Range("AB3").Value = 0.25
$ws.Range("AB4").Value = 0.25

$ws.Range("AC2").Value = 1
$ws.Range("AC3").Value = 1
$ws.Range("AC4").Value = 1

# --- 4. Re-home the two "currently not used" cell comments that used to sit
#        on the (now shifted) dehydrogenation reactor-energy columns so they
#        keep annotating the same logical column (old AJ1/AK1 -> new AL1/AM1)
$commentOld1 = $ws.Range("AJ1").Comment
$commentText1 = $commentOld1.Text()
$ws.Range("AJ1").ClearComments()
$ws.Range("AL1").AddComment($commentText1)

$commentOld2 = $ws.Range("AK1").Comment
$commentText2 = $commentOld2.Text()
$ws.Range("AK1").ClearComments()
$ws.Range("AM1").AddComment($commentText2)

# --- 5. Extend the conditional formatting range so it covers the new columns
$cf = $ws.Range("C3:AL4").FormatConditions
$rule = $cf.Item(1)
$rule.ModifyAppliesToRange($ws.Range("C3:AN4"))

# --- 6. Update the visible selection to reflect the newly added column -----
$ws.Range("AB1").Select()
